$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The SSO block (rows 15-20: "SSO" / "server" / "app_id" / "app_secret" /
# "google client id" / "client secret") is being replaced by a shorter
# "Authentication" block (rows 14-16: "Authentication" / "Authentication
# provider" + "Google / OpenIDConnect / PAM / LDAP" / "credentials"),
# reflecting the deprecation of the SSO server. Delete the now-unneeded
# rows so everything below shifts up to close the gap, then update the
# remaining labels in place.

# Remove the blank row above the "SSO" header so the header itself moves
# from row 15 up to row 14.
$ws.Rows("13:13").Delete()

# Remove three of the now-redundant credential rows ("app_id",
# "app_secret", "google client id") so the block collapses from six rows
# down to three (header + 2 detail rows).
$ws.Rows("16:16").Delete()
$ws.Rows("16:16").Delete()
$ws.Rows("16:16").Delete()

# Relabel the remaining rows of the block.
$ws.Range("A14").Value = "Authentication"
$ws.Range("A15").Value = "Authentication provider"
$ws.Range("B15").Value = "Google / OpenIDConnect / PAM / LDAP"
$ws.Range("A16").Value = "credentials"
